$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.06437833333333333
$ws.Range("H2").Value = 0.193135
$ws.Range("I2").Value = 0.109187438766332
$ws.Range("J2").Value = 0.109187438766332
$ws.Range("M2").Value = 12.39940933333333
$ws.Range("N2").Value = 37.198228
$ws.Range("O2").Value = 0.6889801160127385
$ws.Range("P2").Value = 0.6889801160127385
$ws.Range("Q2").Value = 0.7982533071977777
$ws.Range("R2").Value = 7.18427976478
$ws.Range("S2").Value = 0.07522797422836122
$ws.Range("T2").Value = 0.07522797422836122

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.06437833333333333
$ws.Range("H3").Value = 0.193135
$ws.Range("I3").Value = 0.109187438766332
$ws.Range("J3").Value = 0.109187438766332
$ws.Range("O3").Value = 0.1923952864254561
$ws.Range("P3").Value = 0.1923952864254561
$ws.Range("Q3").Value = 0.2229094426805555
$ws.Range("R3").Value = 2.006184984125
$ws.Range("S3").Value = 0.0210071485555104
$ws.Range("T3").Value = 0.0210071485555104

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.06437833333333333
$ws.Range("H4").Value = 0.193135
$ws.Range("I4").Value = 0.109187438766332
$ws.Range("J4").Value = 0.109187438766332
$ws.Range("M4").Value = 2.134858333333333
$ws.Range("N4").Value = 6.404574999999999
$ws.Range("O4").Value = 0.1186245975618055
$ws.Range("P4").Value = 0.1186245975618055
$ws.Range("Q4").Value = 0.1374386214027778
$ws.Range("R4").Value = 1.236947592625
$ws.Range("S4").Value = 0.01295231598246042
$ws.Range("T4").Value = 0.01295231598246041

$ws.Range("G5").Value = 0.4788196666666666
$ws.Range("I5").Value = 0.8120914339857952
$ws.Range("J5").Value = 0.8120914339857951
$ws.Range("M5").Value = 12.39940933333333
$ws.Range("N5").Value = 37.198228
$ws.Range("O5").Value = 0.6889801160127385
$ws.Range("P5").Value = 0.6889801160127385
$ws.Range("Q5").Value = 5.937081043850221
$ws.Range("R5").Value = 53.433729394652
$ws.Range("S5").Value = 0.5595148504004843
$ws.Range("T5").Value = 0.5595148504004842

$ws.Range("G6").Value = 0.4788196666666666
$ws.Range("I6").Value = 0.8120914339857952
$ws.Range("J6").Value = 0.8120914339857951
$ws.Range("O6").Value = 0.1923952864254561
$ws.Range("P6").Value = 0.1923952864254561
$ws.Range("Q6").Value = 1.657909105669444
$ws.Range("S6").Value = 0.1562425640453564
$ws.Range("T6").Value = 0.1562425640453564

$ws.Range("G7").Value = 0.4788196666666666
$ws.Range("I7").Value = 0.8120914339857952
$ws.Range("J7").Value = 0.8120914339857951
$ws.Range("M7").Value = 2.134858333333333
$ws.Range("N7").Value = 6.404574999999999
$ws.Range("O7").Value = 0.1186245975618055
$ws.Range("P7").Value = 0.1186245975618055
$ws.Range("Q7").Value = 1.022212155547222
$ws.Range("R7").Value = 9.199909399924998
$ws.Range("S7").Value = 0.09633401953995448
$ws.Range("T7").Value = 0.09633401953995445

$ws.Range("G8").Value = 0.042481
$ws.Range("H8").Value = 0.127443
$ws.Range("I8").Value = 0.07204895414449818
$ws.Range("J8").Value = 0.07204895414449818
$ws.Range("M8").Value = 12.39940933333333
$ws.Range("N8").Value = 37.198228
$ws.Range("O8").Value = 0.6889801160127385
$ws.Range("P8").Value = 0.6889801160127385
$ws.Range("Q8").Value = 0.5267393078893333
$ws.Range("R8").Value = 4.740653771004
$ws.Range("S8").Value = 0.04964029678507283
$ws.Range("T8").Value = 0.04964029678507283

$ws.Range("G9").Value = 0.042481
$ws.Range("H9").Value = 0.127443
$ws.Range("I9").Value = 0.07204895414449818
$ws.Range("J9").Value = 0.07204895414449818
$ws.Range("O9").Value = 0.1923952864254561
$ws.Range("P9").Value = 0.1923952864254561
$ws.Range("Q9").Value = 0.1470901084916666
$ws.Range("R9").Value = 1.323810976425
$ws.Range("S9").Value = 0.01386187916928528
$ws.Range("T9").Value = 0.01386187916928528

$ws.Range("G10").Value = 0.042481
$ws.Range("H10").Value = 0.127443
$ws.Range("I10").Value = 0.07204895414449818
$ws.Range("J10").Value = 0.07204895414449818
$ws.Range("M10").Value = 2.134858333333333
$ws.Range("N10").Value = 6.404574999999999
$ws.Range("O10").Value = 0.1186245975618055
$ws.Range("P10").Value = 0.1186245975618055
$ws.Range("Q10").Value = 0.09069091685833333
$ws.Range("R10").Value = 0.816218251725
$ws.Range("S10").Value = 0.008546778190140073
$ws.Range("T10").Value = 0.008546778190140072

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.003934
$ws.Range("H11").Value = 0.011802
$ws.Range("I11").Value = 0.006672173103374587
$ws.Range("J11").Value = 0.006672173103374586
$ws.Range("M11").Value = 12.39940933333333
$ws.Range("N11").Value = 37.198228
$ws.Range("O11").Value = 0.6889801160127385
$ws.Range("P11").Value = 0.6889801160127385
$ws.Range("Q11").Value = 0.04877927631733333
$ws.Range("R11").Value = 0.439013486856
$ws.Range("S11").Value = 0.004596994598820096
$ws.Range("T11").Value = 0.004596994598820095

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.003934
$ws.Range("H12").Value = 0.011802
$ws.Range("I12").Value = 0.006672173103374587
$ws.Range("J12").Value = 0.006672173103374586
$ws.Range("O12").Value = 0.1923952864254561
$ws.Range("P12").Value = 0.1923952864254561
$ws.Range("Q12").Value = 0.01362144221666667
$ws.Range("R12").Value = 0.12259297995
$ws.Range("S12").Value = 0.001283694655303978
$ws.Range("T12").Value = 0.001283694655303977

$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.003934
$ws.Range("H13").Value = 0.011802
$ws.Range("I13").Value = 0.006672173103374587
$ws.Range("J13").Value = 0.006672173103374586
$ws.Range("M13").Value = 2.134858333333333
$ws.Range("N13").Value = 6.404574999999999
$ws.Range("O13").Value = 0.1186245975618055
$ws.Range("P13").Value = 0.1186245975618055
$ws.Range("Q13").Value = 0.008398532683333333
$ws.Range("R13").Value = 0.07558679415
$ws.Range("S13").Value = 0.0007914838492505131
$ws.Range("T13").Value = 0.0007914838492505128

